$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 2.38
$ws.Range("G3").Value = 2.82
$ws.Range("H3").Value = 3.25
$ws.Range("I3").Value = 3.6
$ws.Range("J3").Value = 2.98
$ws.Range("Q3").Value = 2.66
$ws.Range("Q5").Value = 1.36
$ws.Range("S5").Value = 1.89
$ws.Range("X5").Value = 980
$ws.Range("AC5").Value = 980
$ws.Range("AH5").Value = 980
$ws.Range("AL5").Value = 980
$ws.Range("F6").Value = 2.26
$ws.Range("I6").Value = 3.55
$ws.Range("K6").Value = 5.8
$ws.Range("N7").Value = 3.35
$ws.Range("O7").Value = 1.39
$ws.Range("P7").Value = 1.79
$ws.Range("Q7").Value = 2.18
$ws.Range("R7").Value = 1.29
$ws.Range("T7").Value = 1.96
$ws.Range("U7").Value = 1.96
$ws.Range("X7").Value = 13.5
$ws.Range("AB7").Value = 13.5
$ws.Range("AH7").Value = 32
$ws.Range("AI7").Value = 150
$ws.Range("AK7").Value = 400
$ws.Range("G8").Value = 1.22
$ws.Range("I8").Value = 24
$ws.Range("J8").Value = 8
$ws.Range("Q8").Value = 1.52
$ws.Range("F9").Value = 3.05
$ws.Range("N9").Value = 4.5
$ws.Range("O9").Value = 1.27
$ws.Range("Q9").Value = 1.8
$ws.Range("T9").Value = 1.68
$ws.Range("U9").Value = 2.38
$ws.Range("F10").Value = 3.7
$ws.Range("G11").Value = 3.4
$ws.Range("H11").Value = 2.24
$ws.Range("I11").Value = 2.34
$ws.Range("O11").Value = 1.24
$ws.Range("P11").Value = 2.2
$ws.Range("Q11").Value = 1.72
$ws.Range("R11").Value = 1.5
$ws.Range("F12").Value = 1.56
$ws.Range("G12").Value = 1.6
$ws.Range("I12").Value = 6.8
$ws.Range("K12").Value = 4.8
$ws.Range("R12").Value = 1.48
$ws.Range("S12").Value = 2.92
$ws.Range("T12").Value = 1.88
$ws.Range("U12").Value = 2.04
$ws.Range("X12").Value = 26
$ws.Range("Y12").Value = 32
$ws.Range("Z12").Value = 220
$ws.Range("AD12").Value = 36
$ws.Range("AK12").Value = 16
$ws.Range("F13").Value = 2.72
$ws.Range("G13").Value = 2.76
$ws.Range("H13").Value = 2.7
$ws.Range("I13").Value = 2.78
$ws.Range("J13").Value = 3.6
$ws.Range("S13").Value = 2.92
$ws.Range("U13").Value = 2.36
$ws.Range("X13").Value = 22
$ws.Range("Y13").Value = 13.5
$ws.Range("AA13").Value = 980
$ws.Range("AK13").Value = 42
$ws.Range("AL13").Value = 85
$ws.Range("AN13").Value = 20
$ws.Range("AO13").Value = 21
$ws.Range("S14").Value = 2.62
$ws.Range("U14").Value = 2.36
$ws.Range("I15").Value = 7.6
$ws.Range("K15").Value = 5.8
$ws.Range("Q15").Value = 1.5
$ws.Range("R15").Value = 1.73
$ws.Range("G16").Value = 5
$ws.Range("I16").Value = 2.58
$ws.Range("J16").Value = 2.78
$ws.Range("Q16").Value = 3.15
$ws.Range("K17").Value = 4.3
$ws.Range("N17").Value = 3.85
$ws.Range("O17").Value = 1.32
$ws.Range("P17").Value = 1.94
$ws.Range("T17").Value = 1.91
$ws.Range("X17").Value = 15.5
$ws.Range("AB17").Value = 8.6
$ws.Range("AC17").Value = 9.4
$ws.Range("G18").Value = 1.99
$ws.Range("H18").Value = 4.2
$ws.Range("I18").Value = 4.7
$ws.Range("K18").Value = 3.95
$ws.Range("S18").Value = 3.05
$ws.Range("AC18").Value = 9
$ws.Range("F19").Value = 2.52
$ws.Range("G19").Value = 2.54
$ws.Range("J19").Value = 3.5
$ws.Range("K19").Value = 3.6
$ws.Range("O19").Value = 1.31
$ws.Range("P19").Value = 2.02
$ws.Range("F20").Value = 2.82
$ws.Range("G20").Value = 2.9
$ws.Range("H20").Value = 2.62
$ws.Range("I20").Value = 2.66
$ws.Range("P20").Value = 2.02
$ws.Range("S20").Value = 3.1
$ws.Range("T20").Value = 1.72
$ws.Range("U20").Value = 2.28
$ws.Range("X20").Value = 16.5
$ws.Range("F22").Value = 1.87
$ws.Range("G22").Value = 1.96
$ws.Range("H22").Value = 4.1
$ws.Range("K22").Value = 4.3
$ws.Range("P22").Value = 2.28
$ws.Range("F23").Value = 2.02
$ws.Range("F24").Value = 2.8
$ws.Range("G24").Value = 2.86
$ws.Range("H24").Value = 2.72
$ws.Range("I24").Value = 2.78
$ws.Range("Q24").Value = 1.96
$ws.Range("AB24").Value = 12
$ws.Range("AF24").Value = 20
$ws.Range("AG24").Value = 13
$ws.Range("AH24").Value = 17.5
$ws.Range("AL24").Value = 980
$ws.Range("AM24").Value = 110
$ws.Range("AN24").Value = 29
$ws.Range("H25").Value = 2.8
$ws.Range("I25").Value = 2.94
$ws.Range("R25").Value = 1.59
$ws.Range("S25").Value = 2.5
